# Commit: "I made changes in signin page and tree page."
#
# This adds a new "SignIn" worksheet (after the existing "registration"
# sheet) containing a tiny username/password table, wires up a hyperlink
# on the password cell (Excel auto-links "...@..." text typed into a
# cell), and leaves "registration" as the active/selected sheet - matching
# the workbook-level state captured in the target OOXML (activeTab/
# tabSelected).

$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the sheet collection (after "registration").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$signIn = $wb.Worksheets.Add($null, $lastSheet)
$signIn.Name = "SignIn"

# Header row.
$signIn.Range("A1").Value = "username"
$signIn.Range("B1").Value = "password"

# Data row - password written first so shared-string allocation order
# matches (password text registered before the username text).
$signIn.Range("B2").Value = "Dsalgo@1"
$signIn.Range("A2").Value = "sonali"

# The password looks like an email/login ("Dsalgo@1"), which Excel's
# AutoFormat-As-You-Type turns into a mailto hyperlink on the cell.
$signIn.Hyperlinks.Add($signIn.Range("B2"), "mailto:Dsalgo@1")

# Leave the selection on A2 of the new sheet (matches saved view state).
$signIn.Range("A2").Select()

# "registration" ends up the active tab.
$wb.Worksheets.Item("registration").Activate()
